$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell contents (keeps cell formatting/styles intact) so the
# shared-string table can be rebuilt from scratch in the exact order we need.
$ws.Cells.ClearContents()

# Row 1: HKL index header (0..14) - unchanged from before
$row1 = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14)
for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value2 = $row1[$i]
}

# Column B labels for rows 2..19, written in this exact order so that the
# workbook's shared-string table is built in the same order the target file
# expects (the engine assigns shared-string indices in first-seen order).
$labels = @(
    "HKL",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "Gaussian-Quadrature",
    "Spiral-90deg-10rot-5space",
    "Spiral-90deg-15rot-5space",
    "Spiral-90deg-10rot-3space",
    "NoRotation-tilt60deg",
    "Rotation-NoTilt",
    "Rotation-60detTilt",
    "HexGrid-90degTilt5degRes",
    "HexGrid-90degTilt22p5degRes",
    "HexGrid-60degTilt5degRes"
)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value2 = $labels[$i]
}

# Column A index values for rows 2..19 (0-based row index, mirrors column B)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $i
}

# Row 2 (HKL header) columns C..P - written after the B-column labels so
# these strings land at the end of the shared-string table, matching the
# target workbook.
$headers = @(
    "[1, 1, 1]",
    "[2, 0, 0]",
    "[2, 2, 0]",
    "[3, 1, 1]",
    "[2, 2, 2]",
    "[4, 0, 0]",
    "[3, 3, 1]",
    "[4, 2, 0]",
    "[4, 2, 2]",
    "[5, 1, 1]",
    "[3, 3, 3]",
    "2Pairs",
    "4Pairs",
    "MaxUnique"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(2, $i + 3).Value2 = $headers[$i]
}

# Body values: rows 3..19, columns C..P are all 1
for ($r = 3; $r -le 19; $r++) {
    for ($c = 3; $c -le 16; $c++) {
        $ws.Cells.Item($r, $c).Value2 = 1
    }
}
